$d = $word.ActiveDocument

# 1. Cover page document code: C1.008 -> C2.008
[void]$d.Content.Find.Execute("C1.008", $false, $false, $false, $false, $false, $true, 1, $false, "C2.008", 2)

# 2. Remove the "David Guillén Fernández" contributor line entirely (whole paragraph, incl. mark)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*David Guillén Fernández*") {
        $p.Range.Delete()
        break
    }
}

# 3. Update the report date shown under the contributor/team list: 12/02/2024 -> 04/07/2024
#    (only the standalone paragraph right after the team emails; the revision-table date stays)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "12/02/2024`r") {
        $rng = $p.Range
        [void]$rng.MoveEnd(1, -1)
        $rng.Text = "04/07/2024"
        break
    }
}

# 4. Tidy the "Versión inicial del analysis report" revision-table cell: collapse the
#    spell-check-split runs back into plain text (drops stray w:proofErr wrappers).
$cell = $d.Tables(1).Rows(2).Cells(3)
[void]$cell.Range.Find.Execute("Versión inicial del analysis report", $false, $false, $false, $false, $false, $true, 1, $false, "Versión inicial del analysis report", 2)
